$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 5) down into the two
# new rows so the new cells pick up the same styles (centered text / date
# number format) used throughout the table.
$ws.Range("A5:F5").Copy()
$ws.Range("A6:F6").PasteSpecial(-4122)
$ws.Range("A5:F5").Copy()
$ws.Range("A7:F7").PasteSpecial(-4122)

# Row 6: 01_Kirana / Sugar
$ws.Range("A6").Value = "01_Kirana"
$ws.Range("B6").Value = "Sugar"
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 125
$ws.Range("E6").Value = (Get-Date -Year 2023 -Month 1 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F6").Value = (Get-Date -Year 2023 -Month 3 -Day 23 -Hour 0 -Minute 0 -Second 0)

# Row 7: 01_Kirana / Vinegar
$ws.Range("A7").Value = "01_Kirana"
$ws.Range("B7").Value = "Vinegar"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 30
$ws.Range("E7").Value = (Get-Date -Year 2023 -Month 1 -Day 20 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F7").Value = (Get-Date -Year 2023 -Month 1 -Day 20 -Hour 0 -Minute 0 -Second 0)

$ws.Range("E5").Select()
